$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Edit 1: "accellration" -> "accelleration" (insert "e" after "accell")
#         The Word "_GoBack" bookmark ends up right after the inserted
#         "e" (i.e. between "e" and "ration"), exactly where the last
#         edit was made.
# -----------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("accellration")
$insertPos = $rng.Start + 6   # position right after "accell", before "ration"

# Insert the missing "e"
$insRange = $d.Range($insertPos, $insertPos)
$insRange.InsertAfter("e")

# Force the run to split off the newly-typed "e" from the preceding
# text by bookmarking the boundary right before it ...
$bSplitBefore = $d.Range($insertPos, $insertPos)
$d.Bookmarks.Add("ZZZTempSplit", $bSplitBefore)

# ... then place (or move) the document's single "_GoBack" bookmark at
# the boundary right after the "e" - this both matches Word's real
# behaviour (the _GoBack bookmark always tracks the most recent edit)
# and forces the "e" to live in its own run, split from "ration...".
$bAfterE = $d.Range($insertPos + 1, $insertPos + 1)
$d.Bookmarks.Add("_GoBack", $bAfterE)

# Remove the temporary helper bookmark - the "e" stays in its own run.
$d.Bookmarks("ZZZTempSplit").Delete()

# -----------------------------------------------------------------
# Edit 2: "...small structure makes it possible for the business..."
#         -> "...small structure enables the business..."
# -----------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("small structure makes it possible for the business")
$foundText = $rng2.Text
$relOffset = $foundText.IndexOf("makes it possible for")

$split1 = $rng2.Start + 5                      # right after "small"
$replaceStart = $rng2.Start + $relOffset        # start of "makes it possible for"
$replaceEnd = $replaceStart + ("makes it possible for".Length)

# Bookmark both boundaries so the replaced phrase ends up isolated in
# its own run(s), matching the original author's edit pattern.
$b1 = $d.Range($split1, $split1)
$d.Bookmarks.Add("ZZZTempSplit1", $b1)
$b2 = $d.Range($replaceEnd, $replaceEnd)
$d.Bookmarks.Add("ZZZTempSplit2", $b2)

$midRange = $d.Range($replaceStart, $replaceEnd)
$midRange.Text = "enables"

$d.Bookmarks("ZZZTempSplit1").Delete()
$d.Bookmarks("ZZZTempSplit2").Delete()
